# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# per-job Leve tracking sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 8274.857
$ws.Range("I12").Value = 8322.666999999999
$ws.Range("J12").Value = 7988
$ws.Range("K12").Value = 8322.666999999999
$ws.Range("L12").Value = 7988
$ws.Range("M12").Value = -8152.666999999999
$ws.Range("N12").Value = -8328

$ws.Range("H20").Value = 5798
$ws.Range("I20").Value = 5798
$ws.Range("K20").Value = 5798
$ws.Range("M20").Value = -5568

$ws.Range("H34").Value = 2333
$ws.Range("I34").Value = 2333
$ws.Range("K34").Value = 2333
$ws.Range("M34").Value = -2130

$ws.Range("H35").Value = 5798
$ws.Range("I35").Value = 5798
$ws.Range("K35").Value = 5798
$ws.Range("M35").Value = -5419

$ws.Range("H36").Value = 2333
$ws.Range("I36").Value = 2333
$ws.Range("K36").Value = 2333
$ws.Range("M36").Value = -1618

$ws.Range("H38").Value = 334.75
$ws.Range("I38").Value = 239.71428
$ws.Range("K38").Value = 719.14284
$ws.Range("M38").Value = -347.14284

$ws.Range("H48").Value = 20499.5
$ws.Range("I48").Value = 15000
$ws.Range("K48").Value = 45000
$ws.Range("M48").Value = -44708

$ws.Range("H55").Value = 564.75
$ws.Range("J55").Value = 693.75
$ws.Range("L55").Value = 693.75
$ws.Range("N55").Value = -1121.75

$ws.Range("H56").Value = 20499.5
$ws.Range("I56").Value = 15000
$ws.Range("K56").Value = 45000
$ws.Range("M56").Value = -44466

$ws.Range("H80").Value = 988.7
$ws.Range("I80").Value = 459.75
$ws.Range("J80").Value = 1341.3334
$ws.Range("K80").Value = 1379.25
$ws.Range("L80").Value = 4024.0002
$ws.Range("M80").Value = -381.25
$ws.Range("N80").Value = -6020.0002

$ws.Range("H83").Value = 988.7
$ws.Range("I83").Value = 459.75
$ws.Range("J83").Value = 1341.3334
$ws.Range("K83").Value = 4137.75
$ws.Range("L83").Value = 12072.0006
$ws.Range("M83").Value = 854.25
$ws.Range("N83").Value = -22056.0006

$ws.Range("H132").Value = 768.7368
$ws.Range("I132").Value = 792.05554
$ws.Range("J132").Value = 349
$ws.Range("K132").Value = 2376.16662
$ws.Range("L132").Value = 1047
$ws.Range("M132").Value = 153.83338
$ws.Range("N132").Value = -6107

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2667.558
$ws.Range("I32").Value = 2667.558
$ws.Range("K32").Value = 2667.558
$ws.Range("M32").Value = -2380.558

$ws.Range("H95").Value = 80068.336
$ws.Range("J95").Value = 80068.336
$ws.Range("L95").Value = 80068.336
$ws.Range("N95").Value = -85560.336

$ws.Range("H101").Value = 30599.6
$ws.Range("J101").Value = 30599.6
$ws.Range("L101").Value = 30599.6
$ws.Range("N101").Value = -37089.6

$ws.Range("H102").Value = 1400
$ws.Range("I102").Value = 1400
$ws.Range("K102").Value = 1400
$ws.Range("M102").Value = 222

$ws.Range("H110").Value = 6213.6
$ws.Range("I110").Value = 5988.636
$ws.Range("J110").Value = 6832.25
$ws.Range("K110").Value = 5988.636
$ws.Range("L110").Value = 6832.25
$ws.Range("M110").Value = -3943.636
$ws.Range("N110").Value = -10922.25

$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 8008
$ws.Range("I29").Value = 16
$ws.Range("J29").Value = 16000
$ws.Range("K29").Value = 16
$ws.Range("L29").Value = 16000
$ws.Range("M29").Value = 273
$ws.Range("N29").Value = -16578

$ws.Range("H134").Value = 6986.2173
$ws.Range("I134").Value = 6814.7896
$ws.Range("K134").Value = 20444.3688
$ws.Range("M134").Value = -17909.3688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1498
$ws.Range("I16").Value = 1247.6666
$ws.Range("K16").Value = 1247.6666
$ws.Range("M16").Value = -960.6666

$ws.Range("H94").Value = 1494.1666
$ws.Range("I94").Value = 1393.2
$ws.Range("J94").Value = 1999
$ws.Range("K94").Value = 1393.2
$ws.Range("L94").Value = 1999
$ws.Range("M94").Value = -942.2
$ws.Range("N94").Value = -2901

$ws.Range("H113").Value = 1498
$ws.Range("I113").Value = 1247.6666
$ws.Range("K113").Value = 1247.6666
$ws.Range("M113").Value = 922.3334

$ws.Range("H132").Value = 1473.5
$ws.Range("I132").Value = 1398.1666
$ws.Range("K132").Value = 4194.4998
$ws.Range("M132").Value = -1664.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6112611
$ws.Range("I4").Value = 278777
$ws.Range("K4").Value = 836331
$ws.Range("M4").Value = -836219

$ws.Range("H17").Value = 820
$ws.Range("I17").Value = 251
$ws.Range("K17").Value = 753
$ws.Range("M17").Value = -584

$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2831

$ws.Range("H39").Value = 8499.833000000001
$ws.Range("J39").Value = 8499.833000000001
$ws.Range("L39").Value = 25499.499
$ws.Range("N39").Value = -26087.499

$ws.Range("H62").Value = 10716.667
$ws.Range("J62").Value = 10775
$ws.Range("L62").Value = 32325
$ws.Range("N62").Value = -33697

$ws.Range("H65").Value = 10716.667
$ws.Range("J65").Value = 10775
$ws.Range("L65").Value = 96975
$ws.Range("N65").Value = -103839

$ws.Range("H92").Value = 339.66666
$ws.Range("I92").Value = 339.66666
$ws.Range("K92").Value = 1018.99998
$ws.Range("M92").Value = 229.0000200000001

$ws.Range("H98").Value = 1156.6666
$ws.Range("I98").Value = 2000
$ws.Range("J98").Value = 735
$ws.Range("K98").Value = 6000
$ws.Range("L98").Value = 2205
$ws.Range("M98").Value = -4502
$ws.Range("N98").Value = -5201

$ws.Range("H140").Value = 590128.75
$ws.Range("I140").Value = 590128.75
$ws.Range("K140").Value = 1770386.25
$ws.Range("M140").Value = -1765206.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3077.2666
$ws.Range("I7").Value = 3154.2144
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 3154.2144
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -3042.2144
$ws.Range("N7").Value = -2224

$ws.Range("H16").Value = 327.63635
$ws.Range("I16").Value = 339.66666
$ws.Range("J16").Value = 273.5
$ws.Range("K16").Value = 339.66666
$ws.Range("L16").Value = 273.5
$ws.Range("M16").Value = -169.66666
$ws.Range("N16").Value = -613.5

$ws.Range("H61").Value = 1874.5
$ws.Range("J61").Value = 1999
$ws.Range("L61").Value = 1999
$ws.Range("N61").Value = -2403

$ws.Range("H113").Value = 1874.5
$ws.Range("J113").Value = 1999
$ws.Range("L113").Value = 1999
$ws.Range("N113").Value = -6339

$ws.Range("H126").Value = 3077.2666
$ws.Range("I126").Value = 3154.2144
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 9462.643199999999
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -6992.643199999999
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 37305.4
$ws.Range("I45").Value = 24436.428
$ws.Range("J45").Value = 67333
$ws.Range("K45").Value = 24436.428
$ws.Range("L45").Value = 67333
$ws.Range("M45").Value = -23945.428
$ws.Range("N45").Value = -68315
